$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-278.434814453125, -0.2089, 0.0071000000461936, 0.1920000016689301, -1.067000031471252),
    @(17.83259963989258, 0.0143, 0.0353, 0.1861000061035156, -0.5751000046730042),
    @(-2.079600095748901, -0.0024, 0.0186, 0.1064999997615814, -0.3657000064849854),
    @(66.30210113525391, 0.0699, 0.0398, 0.8353999853134155, -0.2746999859809875),
    @(112.2123031616211, 0.0888, 0.033, 1.145699977874756, -0.3203000128269196),
    @(-6.190400123596191, -0.0062, 0.01820000074803829, 0.65420001745224, -0.4032999873161316),
    @(-51.31719970703125, -0.046, -0.0607, 0.65420001745224, -0.4408999979496002),
    @(-141.6749877929688, -0.0186, 0.0257, 1.145699977874756, -1.067000031471252)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
